$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("L2").Value = 6110
$ws.Range("L3").Value = 6644
$ws.Range("K4").Value = 1795
$ws.Range("L4").Value = 1634
$ws.Range("L5").Value = 396
$ws.Range("L6").Value = 5444
$ws.Range("K7").Value = 27587
$ws.Range("L7").Value = 20228
$ws = $wb.Worksheets.Item(2)
$ws.Range("L21").Value = 66
$ws.Range("L22").Value = 65
$ws.Range("L26").Value = 17
$ws.Range("L29").Value = 1133
$ws.Range("L33").Value = 913
$ws.Range("L36").Value = 257
$ws.Range("L37").Value = 775
$ws.Range("L47").Value = 143
$ws.Range("L51").Value = 255
$ws.Range("L53").Value = 223
$ws.Range("L54").Value = 438
$ws.Range("L55").Value = 215
$ws.Range("L63").Value = 60
$ws.Range("L67").Value = 700
$ws.Range("L68").Value = 64
$ws.Range("K71").Value = 83
$ws.Range("L78").Value = 266
$ws.Range("L79").Value = 558
$ws.Range("L83").Value = 444
$ws.Range("L84").Value = 193
$ws.Range("L85").Value = 1003
$ws.Range("L88").Value = 216
$ws.Range("L89").Value = 279
$ws.Range("L90").Value = 213
$ws.Range("L94").Value = 251
$ws.Range("L97").Value = 165
$ws.Range("L99").Value = 350
$ws.Range("K101").Value = 27587
$ws.Range("L101").Value = 20228
$ws = $wb.Worksheets.Item(5)
$ws.Range("L2").Value = 225
$ws.Range("L3").Value = 206
$ws = $wb.Worksheets.Item(7)
$ws.Range("L6").Value = 80
$ws.Range("L7").Value = 279
$ws = $wb.Worksheets.Item(8)
$ws.Range("L2").Value = 303
$ws.Range("L7").Value = 1003
$ws = $wb.Worksheets.Item(11)
$ws.Range("L2").Value = 68
$ws.Range("L7").Value = 223
$ws = $wb.Worksheets.Item(12)
$ws.Range("L2").Value = 402
$ws.Range("L3").Value = 475
$ws = $wb.Worksheets.Item(13)
$ws.Range("L2").Value = 139
$ws.Range("L3").Value = 181
$ws.Range("L7").Value = 444
$ws = $wb.Worksheets.Item(14)
$ws.Range("L6").Value = 262
$ws.Range("L7").Value = 913
$ws = $wb.Worksheets.Item(16)
$ws.Range("L3").Value = 275
$ws.Range("L4").Value = 42
$ws.Range("L6").Value = 203
$ws.Range("L7").Value = 775
$ws = $wb.Worksheets.Item(18)
$ws.Range("L3").Value = 142
$ws.Range("L7").Value = 350
$ws = $wb.Worksheets.Item(21)
$ws.Range("L2").Value = 198
$ws.Range("L7").Value = 700
$ws = $wb.Worksheets.Item(22)
$ws.Range("L2").Value = 63
$ws.Range("L7").Value = 193
$ws = $wb.Worksheets.Item(24)
$ws.Range("L6").Value = 210
$ws.Range("L7").Value = 438
$ws = $wb.Worksheets.Item(25)
$ws.Range("L2").Value = 340
$ws.Range("L3").Value = 437
$ws.Range("L6").Value = 277
$ws.Range("L7").Value = 1133
$ws = $wb.Worksheets.Item(35)
$ws.Range("L4").Value = 31
$ws.Range("L6").Value = 75
$ws.Range("L7").Value = 266
$ws = $wb.Worksheets.Item(36)
$ws.Range("L2").Value = 64
$ws.Range("L7").Value = 215
$ws = $wb.Worksheets.Item(41)
$ws.Range("L6").Value = 33
$ws.Range("L7").Value = 66
$ws = $wb.Worksheets.Item(42)
$ws.Range("L2").Value = 175
$ws.Range("L7").Value = 558
$ws = $wb.Worksheets.Item(44)
$ws.Range("L3").Value = 178
$ws.Range("L4").Value = 49
$ws = $wb.Worksheets.Item(47)
$ws.Range("L5").Value = 3
$ws.Range("L6").Value = 63
$ws.Range("L7").Value = 257
$ws = $wb.Worksheets.Item(51)
$ws.Range("L4").Value = 33
$ws.Range("L7").Value = 251
$ws = $wb.Worksheets.Item(53)
$ws.Range("L2").Value = 54
$ws.Range("L3").Value = 48
$ws.Range("L7").Value = 143
$ws = $wb.Worksheets.Item(57)
$ws.Range("L3").Value = 2
$ws.Range("L7").Value = 17
$ws = $wb.Worksheets.Item(65)
$ws.Range("L3").Value = 35
$ws.Range("L7").Value = 165
$ws = $wb.Worksheets.Item(68)
$ws.Range("L2").Value = 67
$ws.Range("L7").Value = 216
$ws = $wb.Worksheets.Item(74)
$ws.Range("L3").Value = 62
$ws.Range("L7").Value = 213
$ws = $wb.Worksheets.Item(75)
$ws.Range("L5").Value = 6
$ws.Range("L6").Value = 54
$ws.Range("L7").Value = 255
$ws = $wb.Worksheets.Item(76)
$ws.Range("L3").Value = 18
$ws.Range("L7").Value = 64
$ws = $wb.Worksheets.Item(80)
$ws.Range("L2").Value = 22
$ws.Range("L7").Value = 65
$ws = $wb.Worksheets.Item(81)
$ws.Range("K4").Value = 1
$ws.Range("K7").Value = 83
